$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted columns (D = Price) keep their exact textual representation
# by forcing the cell number format to Text ("@") before assigning number-looking strings.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.992.93"
$ws.Range("E2").Value = "  +2.94%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.122.85"
$ws.Range("E3").Value = "  +10.78%  "

# Row 4
$ws.Range("E4").Value = "  -0.22%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.67"
$ws.Range("E5").Value = "  +4.86%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.17%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5353"
$ws.Range("E7").Value = "  +5.58%  "

# Row 8
$ws.Range("E8").Value = "  +8.74%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09073"
$ws.Range("E9").Value = "  +8.96%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.59"
$ws.Range("E10").Value = "  +10.80%  "

# Row 11
$ws.Range("E11").Value = "  +6.26%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.44"
$ws.Range("E12").Value = "  +5.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.120.48"
$ws.Range("E13").Value = "  +10.68%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.818"
$ws.Range("E14").Value = "  +6.33%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.852"
$ws.Range("E15").Value = "  +8.39%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "98.33"
$ws.Range("E16").Value = "  +6.20%  "

# Row 17
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.06%  "

# Row 18
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001142"
$ws.Range("E18").Value = "  +4.20%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06676"
$ws.Range("E19").Value = "  +2.58%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.30"
$ws.Range("E20").Value = "  +4.34%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.09%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.420"
$ws.Range("E22").Value = "  +7.93%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "31.040.22"
$ws.Range("E23").Value = "  +3.08%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.26"
$ws.Range("E24").Value = "  +8.09%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.368.24"
$ws.Range("E25").Value = "  +10.82%  "

# Row 26
$ws.Range("E26").Value = "  +3.29%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.02"
$ws.Range("E27").Value = "  +5.55%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.588"
$ws.Range("E28").Value = "  +14.98%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.83"
$ws.Range("E29").Value = "  +0.87%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.55"
$ws.Range("E30").Value = "  +4.43%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.184"
$ws.Range("E31").Value = "  +4.17%  "

# Row 32
$ws.Range("E32").Value = "  +3.25%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.297"
$ws.Range("E33").Value = "  +6.18%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.015"
$ws.Range("E34").Value = "  +5.85%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.542"
$ws.Range("E35").Value = "  +27.32%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02635"
$ws.Range("E36").Value = "  +7.82%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "13.36"
$ws.Range("E37").Value = "  +16.84%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.581"
$ws.Range("E38").Value = "  +5.53%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.639"
$ws.Range("E39").Value = "  +12.32%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06779"
$ws.Range("E40").Value = "  +5.38%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2294"
$ws.Range("E41").Value = "  +7.15%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6872"
$ws.Range("E42").Value = "  +6.35%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.261"
$ws.Range("E43").Value = "  +4.12%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6495"
$ws.Range("E44").Value = "  +7.55%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.16"
$ws.Range("E45").Value = "  +6.04%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9998"
$ws.Range("E46").Value = "  -0.09%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.269"
$ws.Range("E47").Value = "  +4.09%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.685"
$ws.Range("E48").Value = "  +1.72%  "

# Row 49
$ws.Range("E49").Value = "  +6.62%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.51"
$ws.Range("E50").Value = "  +8.29%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.181"
$ws.Range("E51").Value = "  +3.87%  "
